$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the outcome label in B2 to include the program id
$ws.Range("B2").Value = "earn_rule_outcome_1565159"

# Update the active selection to B3
$ws.Range("B3").Select()
